# Weekly update: insert a new price record as row 35 ("Fruta / hortaliza, semanal").
# This pushes the existing rows 35-147 down to 36-148 (preserving their data and
# formatting) and populates the new row 35 with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 - shifts rows 35..147 down to 36..148.
$ws.Rows.Item(35).Insert()

$ws.Range("A35").Value = 2
$ws.Range("B35").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 45281
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112030
$ws.Range("G35").Value = "Poroto granado"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 240
$ws.Range("K35").Value = 25000
$ws.Range("L35").Value = 27000
$ws.Range("M35").Value = 26000
$ws.Range("N35").Value = "$/caja 15 kilos"
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 1733
$ws.Range("Q35").Value = 15
$ws.Range("R35").Value = "Hortaliza"
